$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1100.2858
$ws.Range("I33").Value = 1100.2858
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1100.2858
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -871.2858000000001
$ws.Range("N33").ClearContents()
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142
$ws.Range("H74").Value = 3945
$ws.Range("I74").Value = 3945
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3945
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3009
$ws.Range("N74").ClearContents()
$ws.Range("H76").Value = 3300
$ws.Range("I76").Value = 3400
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3400
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -3085
$ws.Range("N76").Value = -3830
$ws.Range("H77").Value = 3945
$ws.Range("I77").Value = 3945
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 19725
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15045
$ws.Range("N77").ClearContents()
$ws.Range("H79").Value = 3300
$ws.Range("I79").Value = 3400
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3400
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -2308
$ws.Range("N79").Value = -5384
$ws.Range("H101").Value = 3202
$ws.Range("J101").Value = 5250
$ws.Range("L101").Value = 15750
$ws.Range("N101").Value = -18994
$ws.Range("H111").Value = 2018.5
$ws.Range("I111").Value = 5179
$ws.Range("J111").Value = 1386.4
$ws.Range("K111").Value = 15537
$ws.Range("L111").Value = 4159.200000000001
$ws.Range("M111").Value = -12470
$ws.Range("N111").Value = -10293.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6036706
$ws.Range("I32").Value = 6955320
$ws.Range("K32").Value = 6955320
$ws.Range("M32").Value = -6955033
$ws.Range("H39").Value = 20000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 20000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 20000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -21040
$ws.Range("H45").Value = 2631.8
$ws.Range("I45").Value = 2161.5
$ws.Range("J45").Value = 3169.2856
$ws.Range("K45").Value = 2161.5
$ws.Range("L45").Value = 3169.2856
$ws.Range("M45").Value = -1784.5
$ws.Range("N45").Value = -3923.2856
$ws.Range("H96").Value = 81333.336
$ws.Range("J96").Value = 81333.336
$ws.Range("L96").Value = 81333.336
$ws.Range("N96").Value = -86825.336
$ws.Range("H97").Value = 1160.1666
$ws.Range("I97").Value = 991.9091
$ws.Range("J97").Value = 3011
$ws.Range("K97").Value = 991.9091
$ws.Range("L97").Value = 3011
$ws.Range("M97").Value = -495.9091
$ws.Range("N97").Value = -4003
$ws.Range("H122").Value = 73285.14
$ws.Range("I122").Value = 167982
$ws.Range("J122").Value = 2262.5
$ws.Range("K122").Value = 503946
$ws.Range("L122").Value = 6787.5
$ws.Range("M122").Value = -501496
$ws.Range("N122").Value = -11687.5
$ws.Range("H132").Value = 1835961.6
$ws.Range("I132").Value = 4245.5
$ws.Range("J132").Value = 3850849.5
$ws.Range("K132").Value = 12736.5
$ws.Range("L132").Value = 11552548.5
$ws.Range("M132").Value = -10206.5
$ws.Range("N132").Value = -11557608.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30832
$ws.Range("H99").Value = 1401.1111
$ws.Range("I99").Value = 1228.5714
$ws.Range("J99").Value = 2005
$ws.Range("K99").Value = 1228.5714
$ws.Range("L99").Value = 2005
$ws.Range("M99").Value = 269.4286
$ws.Range("N99").Value = -5001
$ws.Range("H105").Value = 2372.5
$ws.Range("I105").Value = 2372.5
$ws.Range("K105").Value = 2372.5
$ws.Range("M105").Value = -625.5
$ws.Range("H107").Value = 144242.72
$ws.Range("I107").Value = 333999.66
$ws.Range("K107").Value = 333999.66
$ws.Range("M107").Value = -332079.66
$ws.Range("H134").Value = 3400.8076
$ws.Range("I134").Value = 2958.8572
$ws.Range("K134").Value = 8876.571599999999
$ws.Range("M134").Value = -6341.571599999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 9999
$ws.Range("J32").Value = 9999
$ws.Range("L32").Value = 9999
$ws.Range("N32").Value = -10631
$ws.Range("H35").Value = 1000000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H99").Value = 2473.7144
$ws.Range("I99").Value = 2287.2856
$ws.Range("J99").Value = 2520.3215
$ws.Range("K99").Value = 2287.2856
$ws.Range("L99").Value = 2520.3215
$ws.Range("M99").Value = -789.2856000000002
$ws.Range("N99").Value = -5516.3215
$ws.Range("H126").Value = 2473.7144
$ws.Range("I126").Value = 2287.2856
$ws.Range("J126").Value = 2520.3215
$ws.Range("K126").Value = 6861.8568
$ws.Range("L126").Value = 7560.9645
$ws.Range("M126").Value = -4391.8568
$ws.Range("N126").Value = -12500.9645

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 586.03705
$ws.Range("I113").Value = 599.2857
$ws.Range("J113").Value = 571.7692
$ws.Range("K113").Value = 1797.8571
$ws.Range("L113").Value = 1715.3076
$ws.Range("M113").Value = 372.1428999999998
$ws.Range("N113").Value = -6055.3076
$ws.Range("H132").Value = 2601.282
$ws.Range("I132").Value = 2370.6956
$ws.Range("J132").Value = 2932.75
$ws.Range("K132").Value = 21336.2604
$ws.Range("L132").Value = 26394.75
$ws.Range("M132").Value = -18806.2604
$ws.Range("N132").Value = -31454.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9014286
$ws.Range("I11").Value = 12400000
$ws.Range("J11").Value = 550002
$ws.Range("K11").Value = 12400000
$ws.Range("L11").Value = 550002
$ws.Range("M11").Value = -12399861
$ws.Range("N11").Value = -550280
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H70").Value = 5539.478
$ws.Range("I70").Value = 5474.1055
$ws.Range("K70").Value = 5474.1055
$ws.Range("M70").Value = -5204.1055
$ws.Range("H73").Value = 5539.478
$ws.Range("I73").Value = 5474.1055
$ws.Range("K73").Value = 5474.1055
$ws.Range("M73").Value = -4538.1055
$ws.Range("H80").Value = 1881866.6
$ws.Range("J80").Value = 135444.33
$ws.Range("L80").Value = 135444.33
$ws.Range("N80").Value = -137440.33
$ws.Range("H83").Value = 1881866.6
$ws.Range("J83").Value = 135444.33
$ws.Range("L83").Value = 677221.6499999999
$ws.Range("N83").Value = -687205.6499999999
$ws.Range("H102").Value = 1695.1666
$ws.Range("I102").Value = 1813.2941
$ws.Range("J102").Value = 1408.2858
$ws.Range("K102").Value = 1813.2941
$ws.Range("L102").Value = 1408.2858
$ws.Range("M102").Value = -191.2941000000001
$ws.Range("N102").Value = -4652.2858
$ws.Range("H126").Value = 4133.3335
$ws.Range("I126").Value = 3700
$ws.Range("K126").Value = 11100
$ws.Range("M126").Value = -8630
$ws.Range("H132").Value = 45461600
$ws.Range("I132").Value = 52638800
$ws.Range("J132").Value = 6004.3335
$ws.Range("K132").Value = 157916400
$ws.Range("L132").Value = 18013.0005
$ws.Range("M132").Value = -157913870
$ws.Range("N132").Value = -23073.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4792.2
$ws.Range("I7").Value = 4931.9165
$ws.Range("J7").Value = 4233.3335
$ws.Range("K7").Value = 4931.9165
$ws.Range("L7").Value = 4233.3335
$ws.Range("M7").Value = -4819.9165
$ws.Range("N7").Value = -4457.3335
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H40").Value = 3065.1428
$ws.Range("I40").Value = 3076
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3076
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2940
$ws.Range("N40").Value = -3272
$ws.Range("H55").Value = 550
$ws.Range("I55").Value = 444
$ws.Range("J55").Value = 656
$ws.Range("K55").Value = 444
$ws.Range("L55").Value = 656
$ws.Range("M55").Value = -271
$ws.Range("N55").Value = -1002
$ws.Range("H122").Value = 5900.9644
$ws.Range("I122").Value = 4962.6665
$ws.Range("J122").Value = 6013.56
$ws.Range("K122").Value = 14887.9995
$ws.Range("L122").Value = 18040.68
$ws.Range("M122").Value = -12437.9995
$ws.Range("N122").Value = -22940.68
$ws.Range("H126").Value = 4792.2
$ws.Range("I126").Value = 4931.9165
$ws.Range("J126").Value = 4233.3335
$ws.Range("K126").Value = 14795.7495
$ws.Range("L126").Value = 12700.0005
$ws.Range("M126").Value = -12325.7495
$ws.Range("N126").Value = -17640.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 10000
$ws.Range("L56").Value = 10000
$ws.Range("N56").Value = -11428
$ws.Range("H92").Value = 34719.75
$ws.Range("J92").Value = 34719.75
$ws.Range("L92").Value = 34719.75
$ws.Range("N92").Value = -39711.75
$ws.Range("H126").Value = 1495.3334
$ws.Range("I126").Value = 1604.4
$ws.Range("J126").Value = 950
$ws.Range("K126").Value = 4813.200000000001
$ws.Range("L126").Value = 2850
$ws.Range("M126").Value = -2343.200000000001
$ws.Range("N126").Value = -7790
